$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.86240033333333
$ws.Range("H2").Value = 143.587201
$ws.Range("I2").Value = 0.1228118231805696
$ws.Range("J2").Value = 0.1228118231805696
$ws.Range("M2").Value = 4.821747
$ws.Range("N2").Value = 14.465241
$ws.Range("O2").Value = 0.06513881839368539
$ws.Range("P2").Value = 0.0651388183936854
$ws.Range("Q2").Value = 230.780385220049
$ws.Range("R2").Value = 2077.023466980441
$ws.Range("S2").Value = 0.007999817046756527
$ws.Range("T2").Value = 0.007999817046756527
$ws.Range("G3").Value = 47.86240033333333
$ws.Range("H3").Value = 143.587201
$ws.Range("I3").Value = 0.1228118231805696
$ws.Range("J3").Value = 0.1228118231805696
$ws.Range("O3").Value = 0.6589631236492138
$ws.Range("P3").Value = 0.6589631236492139
$ws.Range("Q3").Value = 2334.641113728196
$ws.Range("R3").Value = 21011.77002355376
$ws.Range("S3").Value = 0.08092846262412309
$ws.Range("T3").Value = 0.08092846262412309
$ws.Range("G4").Value = 47.86240033333333
$ws.Range("H4").Value = 143.587201
$ws.Range("I4").Value = 0.1228118231805696
$ws.Range("J4").Value = 0.1228118231805696
$ws.Range("M4").Value = 9.394689
$ws.Range("N4").Value = 28.184067
$ws.Range("O4").Value = 0.1269164351916751
$ws.Range("P4").Value = 0.1269164351916751
$ws.Range("Q4").Value = 449.652365925163
$ws.Range("R4").Value = 4046.871293326466
$ws.Range("S4").Value = 0.01558683879746823
$ws.Range("T4").Value = 0.01558683879746823
$ws.Range("G5").Value = 47.86240033333333
$ws.Range("H5").Value = 143.587201
$ws.Range("I5").Value = 0.1228118231805696
$ws.Range("J5").Value = 0.1228118231805696
$ws.Range("M5").Value = 5.127123333333333
$ws.Range("N5").Value = 15.38137
$ws.Range("O5").Value = 0.0692642636977898
$ws.Range("P5").Value = 0.06926426369778982
$ws.Range("Q5").Value = 245.3964295383745
$ws.Range("R5").Value = 2208.56786584537
$ws.Range("S5").Value = 0.00850647050598531
$ws.Range("T5").Value = 0.00850647050598531
$ws.Range("G6").Value = 47.86240033333333
$ws.Range("H6").Value = 143.587201
$ws.Range("I6").Value = 0.1228118231805696
$ws.Range("J6").Value = 0.1228118231805696
$ws.Range("M6").Value = 5.900888999999999
$ws.Range("N6").Value = 17.702667
$ws.Range("O6").Value = 0.07971735906763582
$ws.Range("P6").Value = 0.07971735906763582
$ws.Range("Q6").Value = 282.4307116405629
$ws.Range("R6").Value = 2541.876404765067
$ws.Range("S6").Value = 0.009790234206236469
$ws.Range("T6").Value = 0.009790234206236467
$ws.Range("I7").Value = 0.04786922362394307
$ws.Range("J7").Value = 0.04786922362394307
$ws.Range("M7").Value = 4.821747
$ws.Range("N7").Value = 14.465241
$ws.Range("O7").Value = 0.06513881839368539
$ws.Range("P7").Value = 0.0651388183936854
$ws.Range("Q7").Value = 89.95288549601202
$ws.Range("R7").Value = 809.5759694641081
$ws.Range("S7").Value = 0.003118144664286742
$ws.Range("T7").Value = 0.003118144664286743
$ws.Range("I8").Value = 0.04786922362394307
$ws.Range("J8").Value = 0.04786922362394307
$ws.Range("O8").Value = 0.6589631236492138
$ws.Range("P8").Value = 0.6589631236492139
$ws.Range("Q8").Value = 909.9894021635854
$ws.Range("R8").Value = 8189.904619472269
$ws.Range("S8").Value = 0.03154405312589627
$ws.Range("T8").Value = 0.03154405312589628
$ws.Range("I9").Value = 0.04786922362394307
$ws.Range("J9").Value = 0.04786922362394307
$ws.Range("M9").Value = 9.394689
$ws.Range("N9").Value = 28.184067
$ws.Range("O9").Value = 0.1269164351916751
$ws.Range("P9").Value = 0.1269164351916751
$ws.Range("Q9").Value = 175.264148842244
$ws.Range("R9").Value = 1577.377339580196
$ws.Range("S9").Value = 0.006075391217743973
$ws.Range("T9").Value = 0.006075391217743974
$ws.Range("I10").Value = 0.04786922362394307
$ws.Range("J10").Value = 0.04786922362394307
$ws.Range("M10").Value = 5.127123333333333
$ws.Range("N10").Value = 15.38137
$ws.Range("O10").Value = 0.0692642636977898
$ws.Range("P10").Value = 0.06926426369778982
$ws.Range("Q10").Value = 95.64988335706224
$ws.Range("R10").Value = 860.84895021356
$ws.Range("S10").Value = 0.003315626528097262
$ws.Range("T10").Value = 0.003315626528097263
$ws.Range("I11").Value = 0.04786922362394307
$ws.Range("J11").Value = 0.04786922362394307
$ws.Range("M11").Value = 5.900888999999999
$ws.Range("N11").Value = 17.702667
$ws.Range("O11").Value = 0.07971735906763582
$ws.Range("P11").Value = 0.07971735906763582
$ws.Range("Q11").Value = 110.084994617444
$ws.Range("R11").Value = 990.7649515569959
$ws.Range("S11").Value = 0.003816008087918825
$ws.Range("T11").Value = 0.003816008087918825
$ws.Range("G12").Value = 171.0598806666667
$ws.Range("H12").Value = 513.1796420000001
$ws.Range("I12").Value = 0.4389285884413335
$ws.Range("J12").Value = 0.4389285884413335
$ws.Range("M12").Value = 4.821747
$ws.Range("N12").Value = 14.465241
$ws.Range("O12").Value = 0.06513881839368539
$ws.Range("P12").Value = 0.0651388183936854
$ws.Range("Q12").Value = 824.8074664248581
$ws.Range("R12").Value = 7423.267197823723
$ws.Range("S12").Value = 0.0285912896102767
$ws.Range("T12").Value = 0.0285912896102767
$ws.Range("G13").Value = 171.0598806666667
$ws.Range("H13").Value = 513.1796420000001
$ws.Range("I13").Value = 0.4389285884413335
$ws.Range("J13").Value = 0.4389285884413335
$ws.Range("O13").Value = 0.6589631236492138
$ws.Range("P13").Value = 0.6589631236492139
$ws.Range("Q13").Value = 8343.990847356352
$ws.Range("R13").Value = 75095.91762620717
$ws.Range("S13").Value = 0.2892377536982413
$ws.Range("T13").Value = 0.2892377536982414
$ws.Range("G14").Value = 171.0598806666667
$ws.Range("H14").Value = 513.1796420000001
$ws.Range("I14").Value = 0.4389285884413335
$ws.Range("J14").Value = 0.4389285884413335
$ws.Range("M14").Value = 9.394689
$ws.Range("N14").Value = 28.184067
$ws.Range("O14").Value = 0.1269164351916751
$ws.Range("P14").Value = 0.1269164351916751
$ws.Range("Q14").Value = 1607.054379240446
$ws.Range("R14").Value = 14463.48941316402
$ws.Range("S14").Value = 0.05570725174868792
$ws.Range("T14").Value = 0.05570725174868794
$ws.Range("G15").Value = 171.0598806666667
$ws.Range("H15").Value = 513.1796420000001
$ws.Range("I15").Value = 0.4389285884413335
$ws.Range("J15").Value = 0.4389285884413335
$ws.Range("M15").Value = 5.127123333333333
$ws.Range("N15").Value = 15.38137
$ws.Range("O15").Value = 0.0692642636977898
$ws.Range("P15").Value = 0.06926426369778982
$ws.Range("Q15").Value = 877.0451055632824
$ws.Range("R15").Value = 7893.405950069541
$ws.Range("S15").Value = 0.03040206549429918
$ws.Range("T15").Value = 0.03040206549429918
$ws.Range("G16").Value = 171.0598806666667
$ws.Range("H16").Value = 513.1796420000001
$ws.Range("I16").Value = 0.4389285884413335
$ws.Range("J16").Value = 0.4389285884413335
$ws.Range("M16").Value = 5.900888999999999
$ws.Range("N16").Value = 17.702667
$ws.Range("O16").Value = 0.07971735906763582
$ws.Range("P16").Value = 0.07971735906763582
$ws.Range("Q16").Value = 1009.405368167246
$ws.Range("R16").Value = 9084.648313505215
$ws.Range("S16").Value = 0.03499022788982832
$ws.Range("T16").Value = 0.03499022788982832
$ws.Range("G17").Value = 12.628047
$ws.Range("H17").Value = 37.884141
$ws.Range("I17").Value = 0.0324027517316099
$ws.Range("J17").Value = 0.0324027517316099
$ws.Range("M17").Value = 4.821747
$ws.Range("N17").Value = 14.465241
$ws.Range("O17").Value = 0.06513881839368539
$ws.Range("P17").Value = 0.0651388183936854
$ws.Range("Q17").Value = 60.88924773810901
$ws.Range("R17").Value = 548.0032296429811
$ws.Range("S17").Value = 0.002110676960501012
$ws.Range("T17").Value = 0.002110676960501013
$ws.Range("G18").Value = 12.628047
$ws.Range("H18").Value = 37.884141
$ws.Range("I18").Value = 0.0324027517316099
$ws.Range("J18").Value = 0.0324027517316099
$ws.Range("O18").Value = 0.6589631236492138
$ws.Range("P18").Value = 0.6589631236492139
$ws.Range("Q18").Value = 615.973238010789
$ws.Range("R18").Value = 5543.759142097101
$ws.Range("S18").Value = 0.02135221849589163
$ws.Range("T18").Value = 0.02135221849589164
$ws.Range("G19").Value = 12.628047
$ws.Range("H19").Value = 37.884141
$ws.Range("I19").Value = 0.0324027517316099
$ws.Range("J19").Value = 0.0324027517316099
$ws.Range("M19").Value = 9.394689
$ws.Range("N19").Value = 28.184067
$ws.Range("O19").Value = 0.1269164351916751
$ws.Range("P19").Value = 0.1269164351916751
$ws.Range("Q19").Value = 118.636574242383
$ws.Range("R19").Value = 1067.729168181447
$ws.Range("S19").Value = 0.004112441740176805
$ws.Range("T19").Value = 0.004112441740176806
$ws.Range("G20").Value = 12.628047
$ws.Range("H20").Value = 37.884141
$ws.Range("I20").Value = 0.0324027517316099
$ws.Range("J20").Value = 0.0324027517316099
$ws.Range("M20").Value = 5.127123333333333
$ws.Range("N20").Value = 15.38137
$ws.Range("O20").Value = 0.0692642636977898
$ws.Range("P20").Value = 0.06926426369778982
$ws.Range("Q20").Value = 64.74555442813001
$ws.Range("R20").Value = 582.70998985317
$ws.Range("S20").Value = 0.002244352740472243
$ws.Range("T20").Value = 0.002244352740472244
$ws.Range("G21").Value = 12.628047
$ws.Range("H21").Value = 37.884141
$ws.Range("I21").Value = 0.0324027517316099
$ws.Range("J21").Value = 0.0324027517316099
$ws.Range("M21").Value = 5.900888999999999
$ws.Range("N21").Value = 17.702667
$ws.Range("O21").Value = 0.07971735906763582
$ws.Range("P21").Value = 0.07971735906763582
$ws.Range("Q21").Value = 74.516703633783
$ws.Range("R21").Value = 670.6503327040469
$ws.Range("S21").Value = 0.002583061794568204
$ws.Range("T21").Value = 0.002583061794568204
$ws.Range("G22").Value = 139.5154473333333
$ws.Range("H22").Value = 418.546342
$ws.Range("I22").Value = 0.3579876130225438
$ws.Range("J22").Value = 0.3579876130225438
$ws.Range("M22").Value = 4.821747
$ws.Range("N22").Value = 14.465241
$ws.Range("O22").Value = 0.06513881839368539
$ws.Range("P22").Value = 0.0651388183936854
$ws.Range("Q22").Value = 672.708189633158
$ws.Range("R22").Value = 6054.373706698422
$ws.Range("S22").Value = 0.0233188901118644
$ws.Range("T22").Value = 0.02331889011186441
$ws.Range("G23").Value = 139.5154473333333
$ws.Range("H23").Value = 418.546342
$ws.Range("I23").Value = 0.3579876130225438
$ws.Range("J23").Value = 0.3579876130225438
$ws.Range("O23").Value = 0.6589631236492138
$ws.Range("P23").Value = 0.6589631236492139
$ws.Range("Q23").Value = 6805.310579413985
$ws.Range("R23").Value = 61247.79521472586
$ws.Range("S23").Value = 0.2359006357050615
$ws.Range("T23").Value = 0.2359006357050615
$ws.Range("G24").Value = 139.5154473333333
$ws.Range("H24").Value = 418.546342
$ws.Range("I24").Value = 0.3579876130225438
$ws.Range("J24").Value = 0.3579876130225438
$ws.Range("M24").Value = 9.394689
$ws.Range("N24").Value = 28.184067
$ws.Range("O24").Value = 0.1269164351916751
$ws.Range("P24").Value = 0.1269164351916751
$ws.Range("Q24").Value = 1310.704238392546
$ws.Range("R24").Value = 11796.33814553291
$ws.Range("S24").Value = 0.04543451168759814
$ws.Range("T24").Value = 0.04543451168759816
$ws.Range("G25").Value = 139.5154473333333
$ws.Range("H25").Value = 418.546342
$ws.Range("I25").Value = 0.3579876130225438
$ws.Range("J25").Value = 0.3579876130225438
$ws.Range("M25").Value = 5.127123333333333
$ws.Range("N25").Value = 15.38137
$ws.Range("O25").Value = 0.0692642636977898
$ws.Range("P25").Value = 0.06926426369778982
$ws.Range("Q25").Value = 715.3129053831711
$ws.Range("R25").Value = 6437.81614844854
$ws.Range("S25").Value = 0.02479574842893581
$ws.Range("T25").Value = 0.02479574842893581
$ws.Range("G26").Value = 139.5154473333333
$ws.Range("H26").Value = 418.546342
$ws.Range("I26").Value = 0.3579876130225438
$ws.Range("J26").Value = 0.3579876130225438
$ws.Range("M26").Value = 5.900888999999999
$ws.Range("N26").Value = 17.702667
$ws.Range("O26").Value = 0.07971735906763582
$ws.Range("P26").Value = 0.07971735906763582
$ws.Range("Q26").Value = 823.2651684993459
$ws.Range("R26").Value = 7409.386516494113
$ws.Range("S26").Value = 0.02853782708908398
$ws.Range("T26").Value = 0.02853782708908398
